# Add a new sales row (row 4) to the "ventas" worksheet, mirroring the
# existing rows 2-3: a date in column A (formatted as YYYY-MM-DD, matching
# the style already used by A2/A3), a product name, a price type, quantity,
# unit value, total value and a discount.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give A4 the same date number format as A2/A3 before assigning its value,
# so it picks up the existing "YYYY-MM-DD" style instead of the default one.
$ws.Range("A4").NumberFormat = "YYYY-MM-DD"
$ws.Range("A4").Value = 45844

$ws.Range("B4").Value = "Poker x30"
$ws.Range("C4").Value = "Unitario"
$ws.Range("D4").Value = 10
$ws.Range("E4").Value = 3000
$ws.Range("F4").Value = 30000
$ws.Range("G4").Value = 0
